# Insert a new row at 306 (this pushes former rows 306..415 down to 307..416,
# Excel automatically updates the sheet dimension to A1:R416).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(306).Insert()

# The row that used to be 306 is now row 307. Duplicate its contents into the
# newly created (empty) row 306, then overwrite the two cells (Fecha / D and
# Volumen / J) that differ for the new record.
$src = $ws.Range("A307:R307")
$dst = $ws.Range("A306:R306")
$src.Copy()
$dst.PasteSpecial(-4104) | Out-Null   # xlPasteAll
$excel.CutCopyMode = 0

$ws.Cells.Item(306, 4).Value = 44900    # D306 Fecha
$ws.Cells.Item(306, 10).Value = 65      # J306 Volumen
